$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 10 data: booking/saleOrder entry
$ws.Range("B10").Value = 445328338
$ws.Range("C10").Value = "KPC"
$ws.Range("D10").Value = "Client"
$ws.Range("E10").Value = "015 855 755/ 085 855 755"

# Match formatting of existing rows: B10:D10 centered like the rest of the table
$ws.Range("B10:D10").HorizontalAlignment = -4108

$ws.Range("E10").Select()
